$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as text (avoids Excel auto-converting
    # numeric-looking strings like "224.01" into real numbers), then
    # restore the "Normal" style so no stray number-format/style id is
    # left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "34.022.73"
Set-TextValue $ws.Range("E2") "  -1.51%  "
Set-TextValue $ws.Range("D3") "1.783.22"
Set-TextValue $ws.Range("E3") "  -3.10%  "
Set-TextValue $ws.Range("E4") "  +0.12%  "
Set-TextValue $ws.Range("D5") "224.01"
Set-TextValue $ws.Range("E5") "  -1.10%  "
Set-TextValue $ws.Range("E6") "  -1.21%  "
Set-TextValue $ws.Range("E7") "  +0.08%  "
Set-TextValue $ws.Range("D8") "32.28"
Set-TextValue $ws.Range("E8") "  -0.52%  "
Set-TextValue $ws.Range("E9") "  -3.65%  "
Set-TextValue $ws.Range("D10") "0.0703"
Set-TextValue $ws.Range("E10") "  -2.03%  "
Set-TextValue $ws.Range("D11") "0.0931"
Set-TextValue $ws.Range("E11") "  -0.34%  "
Set-TextValue $ws.Range("D12") "2.038.73"
Set-TextValue $ws.Range("E12") "  -3.12%  "
Set-TextValue $ws.Range("D13") "1.787.96"
Set-TextValue $ws.Range("E13") "  -2.80%  "
Set-TextValue $ws.Range("D14") "10.79"
Set-TextValue $ws.Range("E14") "  -2.13%  "
Set-TextValue $ws.Range("D15") "34.006.29"
Set-TextValue $ws.Range("E15") "  -1.55%  "
Set-TextValue $ws.Range("D16") "0.621"
Set-TextValue $ws.Range("E16") "  -4.64%  "
Set-TextValue $ws.Range("E17") "  -4.90%  "
Set-TextValue $ws.Range("D18") "67.66"
Set-TextValue $ws.Range("E18") "  -3.18%  "
Set-TextValue $ws.Range("D19") "242.98"
Set-TextValue $ws.Range("E19") "  -3.90%  "
Set-TextValue $ws.Range("E20") "  -3.10%  "
Set-TextValue $ws.Range("E21") "  +0.14%  "
Set-TextValue $ws.Range("D22") "10.64"
Set-TextValue $ws.Range("E22") "  -6.16%  "
Set-TextValue $ws.Range("E23") "  -5.29%  "
Set-TextValue $ws.Range("E24") "  -2.57%  "
Set-TextValue $ws.Range("D25") "159.36"
Set-TextValue $ws.Range("E25") "  -1.47%  "
Set-TextValue $ws.Range("D26") "16.26"
Set-TextValue $ws.Range("E26") "  -3.25%  "
Set-TextValue $ws.Range("E27") "  -3.30%  "
Set-TextValue $ws.Range("E28") "  -2.78%  "
Set-TextValue $ws.Range("E29") "  +0.15%  "
Set-TextValue $ws.Range("E30") "  -4.63%  "
Set-TextValue $ws.Range("E31") "  -0.28%  "
Set-TextValue $ws.Range("E32") "  -4.51%  "
Set-TextValue $ws.Range("E33") "  -4.15%  "
Set-TextValue $ws.Range("E34") "  -7.76%  "
Set-TextValue $ws.Range("D35") "1.391.50"
Set-TextValue $ws.Range("E35") "  -4.70%  "
Set-TextValue $ws.Range("D36") "0.641"
Set-TextValue $ws.Range("E36") "  -2.03%  "
Set-TextValue $ws.Range("E37") "  -3.21%  "
Set-TextValue $ws.Range("E38") "  -4.61%  "
Set-TextValue $ws.Range("E39") "  -1.05%  "
Set-TextValue $ws.Range("D40") "2.20"
Set-TextValue $ws.Range("E40") "  +2.16%  "
Set-TextValue $ws.Range("E41") "  -4.29%  "
Set-TextValue $ws.Range("D42") "78.23"
Set-TextValue $ws.Range("E42") "  -5.76%  "
Set-TextValue $ws.Range("E43") "  -7.65%  "
Set-TextValue $ws.Range("D44") "0.0₆0145"
Set-TextValue $ws.Range("E44") "  +13.41%  "
Set-TextValue $ws.Range("E45") "  +1.22%  "
Set-TextValue $ws.Range("E46") "  -0.07%  "
Set-TextValue $ws.Range("D47") "106.86"
Set-TextValue $ws.Range("E47") "  +0.38%  "
Set-TextValue $ws.Range("D48") "5.85"
Set-TextValue $ws.Range("E48") "  -4.89%  "
Set-TextValue $ws.Range("D49") "12.29"
Set-TextValue $ws.Range("E49") "  -0.38%  "
Set-TextValue $ws.Range("D50") "1.938.92"
Set-TextValue $ws.Range("E50") "  -3.11%  "
Set-TextValue $ws.Range("E51") "  -0.32%  "
